# Insert a new data row at row 25 (weekly price record), shifting the
# existing rows 25-28 down to 26-29, then populate the new row 25 with
# its own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 25; rows 25-28 shift down to 26-29.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new weekly record.
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44809
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = 100112043
$ws.Range("G25").Value = "Pepino dulce"
$ws.Range("H25").Value = "Cultivar IV Región"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = 14000
$ws.Range("L25").Value = 15000
$ws.Range("M25").Value = 14500
$ws.Range("N25").Value = "$/bandeja 18 kilos"
$ws.Range("O25").Value = "Provincia de Limarí"
$ws.Range("P25").Value = 806
$ws.Range("Q25").Value = 18
$ws.Range("R25").Value = "Hortaliza"

# Make sure the D25 date cell uses the same date/time number format as
# the other date cells in column D (style index 2 in the original file).
$ws.Range("D25").NumberFormat = $ws.Range("D26").NumberFormat
